{"js": "// Replace each equation in the 20x5 results table (in document order,\n// row-major) with the updated equation from the commit. Using positional\n// (row, col) addressing rather than text search-and-replace because some\n// old values are duplicated (e.g. \"8+8=16\" appears twice but maps to two\n// different new values), so we must not rely on uniqueness of the old text.\nconst newValues = [\"29+15=44\", \"20-3=17\", \"16-0=16\", \"56+10=66\", \"88-42=46\", \"15+36=51\", \"70-54=16\", \"57+40=97\", \"35-33=2\", \"43+18=61\", \"64+27=91\", \"74-2=72\", \"47-15=32\", \"91-48=43\", \"88-62=26\", \"30+23=53\", \"90-46=44\", \"13+82=95\", \"89-42=47\", \"90+5=95\", \"29+11=40\", \"88-63=25\", \"38-31=7\", \"53+3=56\", \"14+6=20\", \"19-7=12\", \"54+31=85\", \"18-0=18\", \"36-18=18\", \"10+60=70\", \"67+0=67\", \"48-37=11\", \"99-4=95\", \"37-20=17\", \"74+22=96\", \"99-64=35\", \"34+42=76\", \"1+46=47\", \"77+17=94\", \"66-13=53\", \"99-66=33\", \"70-52=18\", \"14+38=52\", \"80-57=23\", \"97-54=43\", \"45-33=12\", \"27+57=84\", \"22+53=75\", \"17-12=5\", \"83-66=17\", \"49+34=83\", \"87-37=50\", \"75-33=42\", \"24-6=18\", \"68-41=27\", \"86-47=39\", \"51+46=97\", \"78-38=40\", \"83-13=70\", \"47+39=86\", \"16+63=79\", \"58+0=58\", \"79-3=76\", \"81-7=74\", \"74-41=33\", \"30-21=9\", \"35-9=26\", \"81-46=35\", \"74+20=94\", \"12+10=22\", \"51+12=63\", \"35+43=78\", \"64-30=34\", \"32+0=32\", \"70-5=65\", \"67-37=30\", \"61-53=8\", \"63-60=3\", \"11+38=49\", \"68+3=71\", \"76-65=11\", \"53-28=25\", \"60-35=25\", \"77+2=79\", \"39+25=64\", \"41-13=28\", \"62+28=90\", \"71+19=90\", \"99-12=87\", \"87-64=23\", \"50-24=26\", \"76-27=49\", \"55+3=58\", \"37+16=53\", \"2+25=27\", \"55+21=76\", \"55-44=11\", \"8-5=3\", \"10-0=10\", \"28+8=36\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst numCols = 5;\nconst numRows = table.rowCount;\n\nlet idx = 0;\nfor (let r = 0; r < numRows; r++) {\n  for (let c = 0; c < numCols; c++) {\n    if (idx >= newValues.length) break;\n    const cell = table.getCell(r, c);\n    cell.value = newValues[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each equation in the 20x5 results table (in document order,\n# row-major) with the updated equation from the commit. We address cells\n# positionally via Table.Cell(row, col) rather than searching for the old\n# text, because some old values are duplicated (e.g. \"8+8=16\" appears twice\n# but maps to two different new values), so matching by text would be\n# ambiguous.\n$newValues = @(\"29+15=44\", \"20-3=17\", \"16-0=16\", \"56+10=66\", \"88-42=46\", \"15+36=51\", \"70-54=16\", \"57+40=97\", \"35-33=2\", \"43+18=61\", \"64+27=91\", \"74-2=72\", \"47-15=32\", \"91-48=43\", \"88-62=26\", \"30+23=53\", \"90-46=44\", \"13+82=95\", \"89-42=47\", \"90+5=95\", \"29+11=40\", \"88-63=25\", \"38-31=7\", \"53+3=56\", \"14+6=20\", \"19-7=12\", \"54+31=85\", \"18-0=18\", \"36-18=18\", \"10+60=70\", \"67+0=67\", \"48-37=11\", \"99-4=95\", \"37-20=17\", \"74+22=96\", \"99-64=35\", \"34+42=76\", \"1+46=47\", \"77+17=94\", \"66-13=53\", \"99-66=33\", \"70-52=18\", \"14+38=52\", \"80-57=23\", \"97-54=43\", \"45-33=12\", \"27+57=84\", \"22+53=75\", \"17-12=5\", \"83-66=17\", \"49+34=83\", \"87-37=50\", \"75-33=42\", \"24-6=18\", \"68-41=27\", \"86-47=39\", \"51+46=97\", \"78-38=40\", \"83-13=70\", \"47+39=86\", \"16+63=79\", \"58+0=58\", \"79-3=76\", \"81-7=74\", \"74-41=33\", \"30-21=9\", \"35-9=26\", \"81-46=35\", \"74+20=94\", \"12+10=22\", \"51+12=63\", \"35+43=78\", \"64-30=34\", \"32+0=32\", \"70-5=65\", \"67-37=30\", \"61-53=8\", \"63-60=3\", \"11+38=49\", \"68+3=71\", \"76-65=11\", \"53-28=25\", \"60-35=25\", \"77+2=79\", \"39+25=64\", \"41-13=28\", \"62+28=90\", \"71+19=90\", \"99-12=87\", \"87-64=23\", \"50-24=26\", \"76-27=49\", \"55+3=58\", \"37+16=53\", \"2+25=27\", \"55+21=76\", \"55-44=11\", \"8-5=3\", \"10-0=10\", \"28+8=36\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$numCols = 5\n$numRows = $t.Rows.Count\n\n$idx = 0\nfor ($r = 1; $r -le $numRows; $r++) {\n    for ($c = 1; $c -le $numCols; $c++) {\n        if ($idx -ge $newValues.Length) { break }\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
